# Update gh-pages output figures for Nanning comic convention info.
# Sheet "展览" (Exhibition) and "全部类型" (All types) both contain the
# same rows; update the "想去人数" (want-to-go count) column F for the
# two affected events.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8592
    $ws.Range("F4").Value = 390
}
